# ACBalanceManager charts: reorder the Xcg/Ycg estimation-method comparison rows
# so that the SFORZA method is listed before TORENBEEK_1982 (rows swapped),
# on the FUSELAGE and WING sheets.

$wb = $excel.ActiveWorkbook

function Swap-Row {
    param(
        $ws,
        [int]$row1,
        [int]$row2
    )

    $aCol = "A"
    $cCol = "C"

    $a1 = $ws.Range("$aCol$row1").Value2
    $c1 = $ws.Range("$cCol$row1").Value2
    $a2 = $ws.Range("$aCol$row2").Value2
    $c2 = $ws.Range("$cCol$row2").Value2

    $ws.Range("$aCol$row1").Value = $a2
    $ws.Range("$cCol$row1").Value = $c2
    $ws.Range("$aCol$row2").Value = $a1
    $ws.Range("$cCol$row2").Value = $c1
}

# FUSELAGE: Xcg ESTIMATION METHOD COMPARISON table, rows 23-24
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
Swap-Row $wsFuselage 23 24

# WING: Xcg ESTIMATION METHOD COMPARISON (rows 23-24) and
#       Ycg ESTIMATION METHOD COMPARISON (rows 27-28)
$wsWing = $wb.Worksheets.Item("WING")
Swap-Row $wsWing 23 24
Swap-Row $wsWing 27 28
